$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set / update cell values ---
$ws.Range("H1").Value2 = "kia-novo.ru_price"
$ws.Range("I1").Value2 = "kia-novo.ru"
$ws.Range("T69").Value2 = 2240104
$ws.Range("H178").Value2 = 2530000
$ws.Range("I178").Value2 = "https://kia-novo.ru/auto/kia/carnival/cuv"
$ws.Range("H179").Value2 = 932300
$ws.Range("I179").Value2 = "https://kia-novo.ru/auto/kia/ceed/hatchback"
$ws.Range("H180").Value2 = 867300
$ws.Range("I180").Value2 = "https://kia-novo.ru/auto/kia/new-ceedsw/ceedswnew"
$ws.Range("H181").Value2 = 1734900
$ws.Range("I181").Value2 = "https://kia-novo.ru/auto/kia/new-cerato/newceratosedan"
$ws.Range("H183").Value2 = 1398000
$ws.Range("I183").Value2 = "https://kia-novo.ru/auto/kia/k5/sedan"
$ws.Range("H186").Value2 = 3772000
$ws.Range("I186").Value2 = "https://kia-novo.ru/auto/kia/mohave-new/suv"
$ws.Range("H187").Value2 = 566900
$ws.Range("I187").Value2 = "https://kia-novo.ru/auto/kia/picanto-new/hatch"
$ws.Range("H189").Value2 = 712000
$ws.Range("I189").Value2 = "https://kia-novo.ru/auto/kia/rio-new/sed"
$ws.Range("H191").Value2 = 889900
$ws.Range("I191").Value2 = "https://kia-novo.ru/auto/kia/rioxline/x-line"
$ws.Range("H192").Value2 = 1021100
$ws.Range("I192").Value2 = "https://kia-novo.ru/auto/kia/seltos/seltos"
$ws.Range("H194").Value2 = 1669900
$ws.Range("I194").Value2 = "https://kia-novo.ru/auto/kia/sorento-new/cuv"
$ws.Range("H197").Value2 = 1032000
$ws.Range("I197").Value2 = "https://kia-novo.ru/auto/kia/new-soul/new"
$ws.Range("H198").Value2 = 977300
$ws.Range("I198").Value2 = "https://kia-novo.ru/auto/kia/new-sportage/sportagenew"
$ws.Range("H200").Value2 = 2156900
$ws.Range("I200").Value2 = "https://kia-novo.ru/auto/kia/stinger-new/sed"
$ws.Range("D208").Value2 = 638940
$ws.Range("E208").Value2 = "https://centorauto-nsk.ru/avto-new/lada/lada_granta_sport_drive_active_liftback/"
$ws.Range("D209").Value2 = 467100
$ws.Range("E209").Value2 = "https://centorauto-nsk.ru/avto-new/lada/granta_hatchback_new/"
$ws.Range("D213").Value2 = 419940
$ws.Range("E213").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D214").Value2 = 620940
$ws.Range("E214").Value2 = "https://centorauto-nsk.ru/avto-new/lada/granta_drive_active/"
$ws.Range("D216").Value2 = 892800
$ws.Range("E216").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D217").Value2 = 868800
$ws.Range("E217").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D218").Value2 = 780600
$ws.Range("E218").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D219").Value2 = 804600
$ws.Range("E219").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D229").Value2 = 575640
$ws.Range("E229").Value2 = "https://centorauto-nsk.ru/avto-new/lada/ladaniva/"
$ws.Range("D232").Value2 = 642720
$ws.Range("E232").Value2 = "https://centorauto-nsk.ru/avto-new/lada/nivaoffroad/"
$ws.Range("D234").Value2 = 517860
$ws.Range("E234").Value2 = "https://centorauto-nsk.ru/avto-new/lada/urban_3d/"
$ws.Range("D235").Value2 = 551400
$ws.Range("E235").Value2 = "https://centorauto-nsk.ru/avto-new/lada/urban_5d/"
$ws.Range("D242").Value2 = 1410540
$ws.Range("E242").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("D247").Value2 = 1305540
$ws.Range("E247").Value2 = "https://lada-novosib.ru/#models"
$ws.Range("T339").Value2 = 4780000

# --- Clear removed cells ---
$ws.Range("X1").ClearContents()
$ws.Range("Y1").ClearContents()
$ws.Range("X178").ClearContents()
$ws.Range("Y178").ClearContents()
$ws.Range("X179").ClearContents()
$ws.Range("Y179").ClearContents()
$ws.Range("X180").ClearContents()
$ws.Range("Y180").ClearContents()
$ws.Range("X181").ClearContents()
$ws.Range("Y181").ClearContents()
$ws.Range("X183").ClearContents()
$ws.Range("Y183").ClearContents()
$ws.Range("X186").ClearContents()
$ws.Range("Y186").ClearContents()
$ws.Range("X187").ClearContents()
$ws.Range("Y187").ClearContents()
$ws.Range("X189").ClearContents()
$ws.Range("Y189").ClearContents()
$ws.Range("X191").ClearContents()
$ws.Range("Y191").ClearContents()
$ws.Range("X192").ClearContents()
$ws.Range("Y192").ClearContents()
$ws.Range("X194").ClearContents()
$ws.Range("Y194").ClearContents()
$ws.Range("X197").ClearContents()
$ws.Range("Y197").ClearContents()
$ws.Range("X198").ClearContents()
$ws.Range("Y198").ClearContents()
$ws.Range("X200").ClearContents()
$ws.Range("Y200").ClearContents()
$ws.Range("H205").ClearContents()
$ws.Range("I205").ClearContents()
$ws.Range("H207").ClearContents()
$ws.Range("I207").ClearContents()
$ws.Range("H208").ClearContents()
$ws.Range("I208").ClearContents()
$ws.Range("H209").ClearContents()
$ws.Range("I209").ClearContents()
$ws.Range("H211").ClearContents()
$ws.Range("I211").ClearContents()
$ws.Range("H213").ClearContents()
$ws.Range("I213").ClearContents()
$ws.Range("H214").ClearContents()
$ws.Range("I214").ClearContents()
$ws.Range("H216").ClearContents()
$ws.Range("I216").ClearContents()
$ws.Range("H217").ClearContents()
$ws.Range("I217").ClearContents()
$ws.Range("H218").ClearContents()
$ws.Range("I218").ClearContents()
$ws.Range("H219").ClearContents()
$ws.Range("I219").ClearContents()
$ws.Range("H220").ClearContents()
$ws.Range("I220").ClearContents()
$ws.Range("H222").ClearContents()
$ws.Range("I222").ClearContents()
$ws.Range("H223").ClearContents()
$ws.Range("I223").ClearContents()
$ws.Range("H229").ClearContents()
$ws.Range("I229").ClearContents()
$ws.Range("H230").ClearContents()
$ws.Range("I230").ClearContents()
$ws.Range("H231").ClearContents()
$ws.Range("I231").ClearContents()
$ws.Range("H232").ClearContents()
$ws.Range("I232").ClearContents()
$ws.Range("H233").ClearContents()
$ws.Range("I233").ClearContents()
$ws.Range("H234").ClearContents()
$ws.Range("I234").ClearContents()
$ws.Range("H235").ClearContents()
$ws.Range("I235").ClearContents()
$ws.Range("H236").ClearContents()
$ws.Range("I236").ClearContents()
$ws.Range("H238").ClearContents()
$ws.Range("I238").ClearContents()
$ws.Range("H239").ClearContents()
$ws.Range("I239").ClearContents()
$ws.Range("H242").ClearContents()
$ws.Range("I242").ClearContents()
$ws.Range("H243").ClearContents()
$ws.Range("I243").ClearContents()
$ws.Range("H244").ClearContents()
$ws.Range("I244").ClearContents()
$ws.Range("H247").ClearContents()
$ws.Range("I247").ClearContents()
$ws.Range("H248").ClearContents()
$ws.Range("I248").ClearContents()
$ws.Range("H249").ClearContents()
$ws.Range("I249").ClearContents()
